# The data rows (2-25) get their "weekly" fields re-permuted: each row's
# Fecha/Calidad/Volumen/Precio-min/Precio-max/Precio-promedio/Unidad/
# Precio-$/Kg/Kg-o-Unidades block is replaced by the block that used to
# live on a different row (row 21 is untouched - it maps to itself).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new-row -> old-row the data block is pulled from
$mapping = @{
    2  = 13
    3  = 12
    4  = 6
    5  = 22
    6  = 18
    7  = 17
    8  = 19
    9  = 2
    10 = 16
    11 = 14
    12 = 25
    13 = 20
    14 = 8
    15 = 4
    16 = 9
    17 = 7
    18 = 24
    19 = 5
    20 = 10
    21 = 21
    22 = 23
    23 = 3
    24 = 15
    25 = 11
}

# Columns that move together as one block: D,I,J,K,L,M,N,P,Q
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot every source row's values BEFORE any writes, so overlapping
# reads/writes in the permutation don't clobber data we still need.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $srcVals[$c]
    }
}
